$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.262.48'
$ws.Range('D3').Value = '3.892.67'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('E4').Value = '  -0.12%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '522.54'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +7.70%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '143.04'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -2.65%  '
$ws.Range('E7').Value = '  -1.75%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.07%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.717'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -2.73%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.173'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +1.79%  '
$ws.Range('E11').Value = '  -4.34%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '41.84'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -3.10%  '
$ws.Range('D13').Value = '4.517.51'
$ws.Range('E13').Value = '  -0.63%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '10.18'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -5.39%  '
$ws.Range('D15').Value = '3.905.55'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.135'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -0.82%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '13.88'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -2.73%  '
$ws.Range('B18').Value = 'Polygon'
$ws.Range('C18').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '1.22'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +6.48%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '19.61'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -3.00%  '
$ws.Range('D20').Value = '69.183.63'
$ws.Range('E20').Value = '  +1.26%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '423.35'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -1.80%  '
$ws.Range('E22').Value = '  -5.91%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '14.11'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -6.36%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '87.85'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -2.29%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '4.00'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +7.31%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '11.47'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -1.34%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '10.46'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -6.84%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '36.09'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -4.41%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '693.26'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -3.43%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '13.05'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -5.20%  '
$ws.Range('E31').Value = '  -4.47%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '2.80'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -4.36%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '67.77'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +11.55%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.435'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +8.10%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '5.89'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -4.13%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '39.86'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -4.81%  '
$ws.Range('D37').Value = '0.0₃0836'
$ws.Range('E37').Value = '  -6.09%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.149'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +2.95%  '
$ws.Range('E39').Value = '  +0.19%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.0478'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -3.21%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '3.04'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +2.13%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '2.75'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -9.03%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '2.95'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -5.98%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '3.32'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.56%  '
$ws.Range('E46').Value = '  -2.08%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '3.02'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +7.46%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '26.85'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +5.82%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0341'
$ws.Range('E49').Value = '  +2.33%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.712.77'
$ws.Range('E50').Value = '  +10.86%  '
$ws.Range('B51').Value = 'LidoDAOToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '3.27'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -4.75%  '
